# Scheduled-runner data refresh: updates market-board-derived profit figures
# (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit*) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit worksheets.
$wb = $excel.ActiveWorkbook

# Worksheet name -> cell reference -> new value.
$updates = @{
    "ALC" = @{
        "H40" = 1200
        "I40" = 1200
        "J40" = 1200
        "K40" = 1200
        "L40" = 1200
        "M40" = -1025
        "N40" = -1550
        "H125" = 1163866.5
        "J125" = 1357577.5
        "L125" = 12218197.5
        "N125" = -12223117.5
        "H138" = 2758.4792
        "I138" = 2545.3872
        "J138" = 3147.0588
        "K138" = 7636.1616
        "L138" = 9441.1764
        "M138" = -2496.1616
        "N138" = -19721.1764
        "H141" = 3983.8462
        "I141" = 2460.625
        "J141" = 6421
        "K141" = 7381.875
        "L141" = 19263
        "M141" = -2201.875
        "N141" = -29623
    }
    "ARM" = @{
        "H24" = 21900
        "J24" = 21900
        "L24" = 21900
        "N24" = -22648
        "H94" = 32916.668
        "J94" = 32916.668
        "L94" = 32916.668
        "N94" = -34718.668
        "H100" = 21900
        "J100" = 21900
        "L100" = 21900
        "N100" = -24064
    }
    "BSM" = @{
        "H94" = 2209.8125
        "I94" = 1613.0834
        "J94" = 4000
        "K94" = 1613.0834
        "L94" = 4000
        "M94" = -1162.0834
        "N94" = -4902
        "H100" = 29999.8
        "J100" = 29999.8
        "L100" = 29999.8
        "N100" = -32163.8
        "H103" = 20000
        "J103" = 20000
        "L103" = 20000
        "N103" = -22344
    }
    "CRP" = @{
        "H97" = 30167.846
        "J97" = 30167.846
        "L97" = 30167.846
        "N97" = -32149.846
        "H99" = 16688630
        "I99" = 23821812
        "K99" = 23821812
        "M99" = -23820314
        "H126" = 16688630
        "I126" = 23821812
        "K126" = 71465436
        "M126" = -71462966
        "H127" = 54980
        "J127" = 54980
        "L127" = 54980
        "N127" = -64900
    }
    "CUL" = @{
        "H22" = 1200
        "J22" = 1200
        "L22" = 3600
        "N22" = -3938
        "H27" = 1200
        "J27" = 1200
        "L27" = 3600
        "N27" = -3804
        "H40" = 1314.6
        "I40" = 627.1429
        "J40" = 1684.7693
        "K40" = 2508.5716
        "L40" = 6739.0772
        "M40" = -2439.5716
        "N40" = -6877.0772
        "H68" = 2781.923
        "J68" = 3851.9768
        "L68" = 11555.9304
        "N68" = -13177.9304
        "H71" = 2781.923
        "J71" = 3851.9768
        "L71" = 34667.7912
        "N71" = -42779.7912
        "H107" = 466743.38
        "I107" = 732789.5
        "J107" = 1162.65
        "K107" = 2198368.5
        "L107" = 3487.95
        "M107" = -2196448.5
        "N107" = -7327.950000000001
    }
    "GSM" = @{
        "H70" = 2935214.2
        "I70" = 1394223.2
        "J70" = 9099179
        "K70" = 1394223.2
        "L70" = 9099179
        "M70" = -1393953.2
        "N70" = -9099719
        "H73" = 2935214.2
        "I73" = 1394223.2
        "J73" = 9099179
        "K73" = 1394223.2
        "L73" = 9099179
        "M73" = -1393287.2
        "N73" = -9101051
        "H94" = 63244
        "J94" = 63244
        "L94" = 63244
        "N94" = -64596
        "H97" = 11369678
        "I97" = 1228.2354
        "J97" = 50022410
        "K97" = 1228.2354
        "L97" = 50022410
        "M97" = -732.2354
        "N97" = -50023402
        "H100" = 30000
        "J100" = 30000
        "L100" = 30000
        "N100" = -32164
    }
    "LTW" = @{
        "H7" = 1824.75
        "I7" = 1824.75
        "J7" = 0
        "K7" = 1824.75
        "L7" = 0
        "M7" = -1712.75
        "H97" = 0
        "J97" = 0
        "L97" = 0
        "H100" = 5452.8667
        "I100" = 1306.6666
        "J100" = 8217
        "K100" = 1306.6666
        "L100" = 8217
        "M100" = -765.6666
        "N100" = -9299
        "H126" = 1824.75
        "I126" = 1824.75
        "J126" = 0
        "K126" = 5474.25
        "L126" = 0
        "M126" = -3004.25
    }
    "WVR" = @{
        "H94" = 38000
        "I94" = 38000
        "J94" = 0
        "K94" = 38000
        "L94" = 0
        "M94" = -37099
        "H96" = 1920.0869
        "I96" = 1648.5555
        "J96" = 2897.6
        "K96" = 1648.5555
        "L96" = 2897.6
        "M96" = -275.5554999999999
        "N96" = -5643.6
        "H98" = 33666.5
        "J98" = 33666.5
        "L98" = 33666.5
        "N98" = -39656.5
        "H103" = 38000
        "J103" = 38000
        "L103" = 38000
        "N103" = -40344
        "H126" = 8404086
        "I126" = 10204754
        "K126" = 30614262
        "M126" = -30611792
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($cellRef in $cellMap.Keys) {
        $ws.Range($cellRef).Value = $cellMap[$cellRef]
    }
}

# Cells whose cached value is dropped entirely in the refreshed data
# (no LeveProfit figure computed this pass for that column).
$clears = @{
    "LTW" = @("N7", "N97", "N126")
    "WVR" = @("N94")
}

foreach ($sheetName in $clears.Keys) {
    $ws = $wb.Worksheets($sheetName)
    foreach ($cellRef in $clears[$sheetName]) {
        $ws.Range($cellRef).ClearContents()
    }
}
